{"js": "// Replace the date line and each arithmetic-problem table cell with the\n// new values from the commit. Cell pairs are addressed by (row, column)\n// position (NOT a blind global find/replace) because a handful of the\n// original expressions (e.g. \"16+76=92\", \"64+28=92\") repeat verbatim in\n// more than one cell but must resolve to different replacement text, and\n// one cell (\"58+17=75\") is left untouched by the edit.\n\nconst DATE_OLD = \"2025-08-01 Friday\";\nconst DATE_NEW = \"2025-08-02 Saturday\";\n\n// [oldText, newText] for every cell, in row-major order (row 1 col 1,\n// row 1 col 2, ... row 20 col 5). newText === oldText for the one cell\n// the diff leaves unchanged, so that cell is simply skipped.\nconst CELL_PAIRS = [[\"65-57=8\", \"86-58=28\"], [\"71-47=24\", \"70-65=5\"], [\"17+76=93\", \"9+42=51\"], [\"40-22=18\", \"49+48=97\"], [\"22-18=4\", \"80-73=7\"], [\"83-14=69\", \"38+26=64\"], [\"38+56=94\", \"87-78=9\"], [\"26+35=61\", \"50-31=19\"], [\"29+55=84\", \"9+49=58\"], [\"71-69=2\", \"49+26=75\"], [\"91-49=42\", \"29+29=58\"], [\"17+59=76\", \"34+8=42\"], [\"30-28=2\", \"42-19=23\"], [\"65+9=74\", \"9+73=82\"], [\"35-17=18\", \"9+77=86\"], [\"24-17=7\", \"43-7=36\"], [\"29+16=45\", \"49+8=57\"], [\"63-35=28\", \"71-27=44\"], [\"5+19=24\", \"40-19=21\"], [\"62+19=81\", \"78+7=85\"], [\"52-47=5\", \"35-16=19\"], [\"82-24=58\", \"55-46=9\"], [\"56+18=74\", \"82-26=56\"], [\"48+27=75\", \"54-16=38\"], [\"93-58=35\", \"91-84=7\"], [\"16+76=92\", \"95-87=8\"], [\"44+37=81\", \"62-27=35\"], [\"76+9=85\", \"87-79=8\"], [\"64+28=92\", \"14-9=5\"], [\"93-77=16\", \"33-15=18\"], [\"9+13=22\", \"55+6=61\"], [\"5+6=11\", \"81-64=17\"], [\"80-51=29\", \"19+37=56\"], [\"15-8=7\", \"74-27=47\"], [\"70-44=26\", \"9+2=11\"], [\"39+23=62\", \"90-16=74\"], [\"18+64=82\", \"87-68=19\"], [\"65-48=17\", \"26+39=65\"], [\"90-11=79\", \"16+9=25\"], [\"70-25=45\", \"22+39=61\"], [\"65+29=94\", \"45-36=9\"], [\"18+63=81\", \"93-8=85\"], [\"67+19=86\", \"26+27=53\"], [\"83-26=57\", \"41-17=24\"], [\"17+24=41\", \"62-13=49\"], [\"90-73=17\", \"54-28=26\"], [\"35+48=83\", \"33+58=91\"], [\"82-27=55\", \"4+8=12\"], [\"46-7=39\", \"62-45=17\"], [\"9+25=34\", \"80-15=65\"], [\"7+54=61\", \"88+5=93\"], [\"87+6=93\", \"5+56=61\"], [\"71-55=16\", \"81-3=78\"], [\"49+43=92\", \"5+39=44\"], [\"90-43=47\", \"74-35=39\"], [\"50-8=42\", \"92-79=13\"], [\"90-41=49\", \"8+79=87\"], [\"75-57=18\", \"21-9=12\"], [\"63+29=92\", \"95-46=49\"], [\"70-63=7\", \"19+38=57\"], [\"37+18=55\", \"57-29=28\"], [\"67+24=91\", \"45+47=92\"], [\"28+26=54\", \"37+35=72\"], [\"9+17=26\", \"62-54=8\"], [\"8+54=62\", \"35+16=51\"], [\"64+28=92\", \"83-8=75\"], [\"34+28=62\", \"16+69=85\"], [\"9+55=64\", \"74-8=66\"], [\"26+56=82\", \"24+57=81\"], [\"8+49=57\", \"93-8=85\"], [\"73-5=68\", \"64+27=91\"], [\"72-47=25\", \"48-19=29\"], [\"72-27=45\", \"69+8=77\"], [\"4+37=41\", \"57+14=71\"], [\"53-6=47\", \"74+19=93\"], [\"23+19=42\", \"84-28=56\"], [\"58+17=75\", null], [\"29+15=44\", \"61-4=57\"], [\"90-2=88\", \"37+6=43\"], [\"52-25=27\", \"15+49=64\"], [\"93-27=66\", \"83-36=47\"], [\"63-46=17\", \"40-19=21\"], [\"63-24=39\", \"39+2=41\"], [\"63-37=26\", \"9+54=63\"], [\"71-42=29\", \"28+59=87\"], [\"13+18=31\", \"68+25=93\"], [\"91-78=13\", \"49+37=86\"], [\"77+15=92\", \"38+54=92\"], [\"19+8=27\", \"18+23=41\"], [\"81-33=48\", \"46+46=92\"], [\"16+76=92\", \"32-26=6\"], [\"29+53=82\", \"84+7=91\"], [\"9+88=97\", \"84-58=26\"], [\"38+34=72\", \"45+28=73\"], [\"47+7=54\", \"47-28=19\"], [\"34+58=92\", \"20-13=7\"], [\"20-9=11\", \"43-14=29\"], [\"60-51=9\", \"71-56=15\"], [\"75-59=16\", \"28+5=33\"], [\"75-58=17\", \"56+39=95\"]];\n\nasync function replaceInParagraph(paragraph, oldText, newText) {\n  const results = paragraph.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  // Replace only the first match - each cell holds exactly one run with\n  // exactly one occurrence of its expression.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nfunction stripCellMark(text) {\n  // Cell text can carry a trailing cell-mark control character; ignore it\n  // when comparing against the expected plain-text expression.\n  return text.replace(/[\\r\\x07]+$/, \"\");\n}\n\n// 1) Update the date heading above the table.\nconst body = context.document.body;\nconst dateResults = body.paragraphs.getFirstOrNullObject().search(DATE_OLD, { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(DATE_NEW, Word.InsertLocation.replace);\n  await context.sync();\n} else {\n  // Fall back to a body-wide search in case the heading paragraph lookup\n  // above ever returns something unexpected.\n  const fallback = body.search(DATE_OLD, { matchCase: true });\n  fallback.load(\"items\");\n  await context.sync();\n  fallback.items[0].insertText(DATE_NEW, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Walk the answer table row by row, column by column, and swap each\n//    cell's expression for its replacement using the pre-computed plan.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\nconst columnCount = table.values[0].length;\n\nlet cellIndex = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const [oldText, newText] = CELL_PAIRS[cellIndex];\n    cellIndex++;\n    if (newText === null || newText === oldText) {\n      // This is the single cell the commit leaves untouched.\n      continue;\n    }\n    const cell = table.getCell(r, c);\n    cell.load(\"value\");\n    await context.sync();\n    const current = stripCellMark(cell.value);\n    if (current !== oldText) {\n      throw new Error(\n        \"cell (\" + r + \",\" + c + \") was '\" + current + \"', expected '\" + oldText + \"'\"\n      );\n    }\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n    const paragraph = paragraphs.items[0];\n    await replaceInParagraph(paragraph, oldText, newText);\n  }\n}\n", "ps1": "# Replace the date line and each arithmetic-problem table cell with the\n# new values from the commit. Cells are addressed by (row, column)\n# position (NOT a blind Find/Replace across the whole document) because a\n# few of the original expressions (e.g. \"16+76=92\", \"64+28=92\") repeat\n# verbatim in more than one cell but must resolve to different\n# replacement text, and one cell (\"58+17=75\") is left untouched by the\n# edit - matching it positionally avoids corrupting either case.\n\n$d = $word.ActiveDocument\n\n# --- 1) Date heading above the table -----------------------------------\n$dateOld = '2025-08-01 Friday'\n$dateNew = '2025-08-02 Saturday'\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $dateOld\n$find.Replacement.Text = $dateNew\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop - don't wrap, the heading is the first paragraph\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n# MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike, MatchAllWordForms,\n# Forward, Wrap, Format, ReplaceWith, Replace=wdReplaceOne\n$find.Execute($dateOld, $true, $false, $false, $false, $false, $true, 0, $false, $dateNew, 1) | Out-Null\n\n# --- 2) Answer table, cell by cell --------------------------------------\n# [Old, New] pairs in row-major order (row 1 col 1, row 1 col 2, ... row\n# 20 col 5). New equals Old for the single cell the commit leaves\n# untouched, so that cell is simply skipped below.\n$cellPairs = @(\n    @{Old='65-57=8'; New='86-58=28'},\n    @{Old='71-47=24'; New='70-65=5'},\n    @{Old='17+76=93'; New='9+42=51'},\n    @{Old='40-22=18'; New='49+48=97'},\n    @{Old='22-18=4'; New='80-73=7'},\n    @{Old='83-14=69'; New='38+26=64'},\n    @{Old='38+56=94'; New='87-78=9'},\n    @{Old='26+35=61'; New='50-31=19'},\n    @{Old='29+55=84'; New='9+49=58'},\n    @{Old='71-69=2'; New='49+26=75'},\n    @{Old='91-49=42'; New='29+29=58'},\n    @{Old='17+59=76'; New='34+8=42'},\n    @{Old='30-28=2'; New='42-19=23'},\n    @{Old='65+9=74'; New='9+73=82'},\n    @{Old='35-17=18'; New='9+77=86'},\n    @{Old='24-17=7'; New='43-7=36'},\n    @{Old='29+16=45'; New='49+8=57'},\n    @{Old='63-35=28'; New='71-27=44'},\n    @{Old='5+19=24'; New='40-19=21'},\n    @{Old='62+19=81'; New='78+7=85'},\n    @{Old='52-47=5'; New='35-16=19'},\n    @{Old='82-24=58'; New='55-46=9'},\n    @{Old='56+18=74'; New='82-26=56'},\n    @{Old='48+27=75'; New='54-16=38'},\n    @{Old='93-58=35'; New='91-84=7'},\n    @{Old='16+76=92'; New='95-87=8'},\n    @{Old='44+37=81'; New='62-27=35'},\n    @{Old='76+9=85'; New='87-79=8'},\n    @{Old='64+28=92'; New='14-9=5'},\n    @{Old='93-77=16'; New='33-15=18'},\n    @{Old='9+13=22'; New='55+6=61'},\n    @{Old='5+6=11'; New='81-64=17'},\n    @{Old='80-51=29'; New='19+37=56'},\n    @{Old='15-8=7'; New='74-27=47'},\n    @{Old='70-44=26'; New='9+2=11'},\n    @{Old='39+23=62'; New='90-16=74'},\n    @{Old='18+64=82'; New='87-68=19'},\n    @{Old='65-48=17'; New='26+39=65'},\n    @{Old='90-11=79'; New='16+9=25'},\n    @{Old='70-25=45'; New='22+39=61'},\n    @{Old='65+29=94'; New='45-36=9'},\n    @{Old='18+63=81'; New='93-8=85'},\n    @{Old='67+19=86'; New='26+27=53'},\n    @{Old='83-26=57'; New='41-17=24'},\n    @{Old='17+24=41'; New='62-13=49'},\n    @{Old='90-73=17'; New='54-28=26'},\n    @{Old='35+48=83'; New='33+58=91'},\n    @{Old='82-27=55'; New='4+8=12'},\n    @{Old='46-7=39'; New='62-45=17'},\n    @{Old='9+25=34'; New='80-15=65'},\n    @{Old='7+54=61'; New='88+5=93'},\n    @{Old='87+6=93'; New='5+56=61'},\n    @{Old='71-55=16'; New='81-3=78'},\n    @{Old='49+43=92'; New='5+39=44'},\n    @{Old='90-43=47'; New='74-35=39'},\n    @{Old='50-8=42'; New='92-79=13'},\n    @{Old='90-41=49'; New='8+79=87'},\n    @{Old='75-57=18'; New='21-9=12'},\n    @{Old='63+29=92'; New='95-46=49'},\n    @{Old='70-63=7'; New='19+38=57'},\n    @{Old='37+18=55'; New='57-29=28'},\n    @{Old='67+24=91'; New='45+47=92'},\n    @{Old='28+26=54'; New='37+35=72'},\n    @{Old='9+17=26'; New='62-54=8'},\n    @{Old='8+54=62'; New='35+16=51'},\n    @{Old='64+28=92'; New='83-8=75'},\n    @{Old='34+28=62'; New='16+69=85'},\n    @{Old='9+55=64'; New='74-8=66'},\n    @{Old='26+56=82'; New='24+57=81'},\n    @{Old='8+49=57'; New='93-8=85'},\n    @{Old='73-5=68'; New='64+27=91'},\n    @{Old='72-47=25'; New='48-19=29'},\n    @{Old='72-27=45'; New='69+8=77'},\n    @{Old='4+37=41'; New='57+14=71'},\n    @{Old='53-6=47'; New='74+19=93'},\n    @{Old='23+19=42'; New='84-28=56'},\n    @{Old='58+17=75'; New='58+17=75'},\n    @{Old='29+15=44'; New='61-4=57'},\n    @{Old='90-2=88'; New='37+6=43'},\n    @{Old='52-25=27'; New='15+49=64'},\n    @{Old='93-27=66'; New='83-36=47'},\n    @{Old='63-46=17'; New='40-19=21'},\n    @{Old='63-24=39'; New='39+2=41'},\n    @{Old='63-37=26'; New='9+54=63'},\n    @{Old='71-42=29'; New='28+59=87'},\n    @{Old='13+18=31'; New='68+25=93'},\n    @{Old='91-78=13'; New='49+37=86'},\n    @{Old='77+15=92'; New='38+54=92'},\n    @{Old='19+8=27'; New='18+23=41'},\n    @{Old='81-33=48'; New='46+46=92'},\n    @{Old='16+76=92'; New='32-26=6'},\n    @{Old='29+53=82'; New='84+7=91'},\n    @{Old='9+88=97'; New='84-58=26'},\n    @{Old='38+34=72'; New='45+28=73'},\n    @{Old='47+7=54'; New='47-28=19'},\n    @{Old='34+58=92'; New='20-13=7'},\n    @{Old='20-9=11'; New='43-14=29'},\n    @{Old='60-51=9'; New='71-56=15'},\n    @{Old='75-59=16'; New='28+5=33'},\n    @{Old='75-58=17'; New='56+39=95'}\n)\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif ($cellPairs.Count -ne ($rowCount * $colCount)) {\n    Write-Output (\"warning: table is \" + $rowCount + \"x\" + $colCount + \" (\" + ($rowCount*$colCount) + \" cells) but plan has \" + $cellPairs.Count + \" entries\")\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $pair = $cellPairs[$i]\n        $i = $i + 1\n        if ($pair.New -eq $pair.Old) {\n            continue\n        }\n        $cell = $table.Cell($r, $c)\n        # Cell.Range.Text carries trailing cell-mark characters (CR + BEL);\n        # trim them before comparing so the guard below is exact.\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $pair.Old) {\n            Write-Output (\"warning: cell ($r,$c) was '\" + $current + \"', expected '\" + $pair.Old + \"'\")\n        }\n        $cell.Range.Text = $pair.New\n    }\n}\n"}
